$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row above the old blank spacer row (row 9), pushing
#    every row below it down by one. This is the row that will hold the
#    new "Type" field.
# ---------------------------------------------------------------------
$ws.Rows("9:9").Insert(-4121)

# Copy the label / merged-label / value look from the row above (row 2,
# "SO No.") onto the new row 9 cells so the new row matches the same
# template used by the other header fields.
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Merge A9:B9 like the other label rows above it.
$ws.Range("A9:B9").Merge()

# ---------------------------------------------------------------------
# 2. Populate the new row's content.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Type"
$ws.Range("C9").Value = 0

$ws.Range("D9").Value = "(0: for ""Sale"" type, 2: for ""WH transfer"")"
$ws.Range("D9").Font.Name = "Arial"
$ws.Range("D9").Font.Size = 10
$ws.Range("D9").Font.Bold = $false
$ws.Range("D9").Font.Color = 255
$ws.Range("D9").HorizontalAlignment = -4131
$ws.Range("D9").VerticalAlignment = -4160

$ws.Rows("9:9").RowHeight = 20.25

# ---------------------------------------------------------------------
# 3. Update the SO No. value shown in C2 ("SO10095" -> "SO444")
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "SO444"

# ---------------------------------------------------------------------
# 4. Fix up the sort range that used to cover the item rows (it shifted
#    down by one row along with everything else).
# ---------------------------------------------------------------------
$sort = $ws.Sort
$sort.SetRange($ws.Range("B12:H21"))
$sort.Header = -4142
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("B12:B21"))
$sort.Apply()
